{"js": "// Translate the \"Path walking Snail\" facilitator guide table into French.\n// Each edit below targets a single, uniquely-identifying run of text and\n// replaces it in place (Word.InsertLocation.replace) so sibling runs (and\n// their own formatting) are left untouched \u2014 matching the source diff,\n// which only rewrites <w:t> text nodes without touching run properties.\n\nconst body = context.document.body;\n\n// 1) Title row, first run: \"The \" -> \"L'escargot marchant \"\n//    (the bold \"p\" run that follows stays untouched, see edit #2)\nconst titlePrefix = body.search(\"The \", { matchCase: true, matchWholeWord: false });\ntitlePrefix.load(\"text\");\nawait context.sync();\nif (titlePrefix.items.length === 0) {\n  throw new Error('Could not find title run \"The \"');\n}\ntitlePrefix.items[0].insertText(\"L'escargot marchant \", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Title row, third run: \"ath walking Snail\" -> \"ath\"\nconst titleSuffix = body.search(\"ath walking Snail\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nif (titleSuffix.items.length === 0) {\n  throw new Error('Could not find title run \"ath walking Snail\"');\n}\ntitleSuffix.items[0].insertText(\"ath\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Objective(s) cell\nconst objective = body.search(\n  \"Understand how to model a difficult mathematical problem. Learn how to solve a puzzle by means of a coloring of the plane.\",\n  { matchCase: true, matchWholeWord: false }\n);\nawait context.sync();\nif (objective.items.length === 0) {\n  throw new Error(\"Could not find objective text\");\n}\nobjective.items[0].insertText(\n  \" Apprenez \u00e0 r\u00e9soudre une \u00e9nigme \u00e0 l'aide d'un coloriage du plan.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 4) Required resources cell\nconst resources = body.search(\n  \"Blank paper (5 per group of students), pen (2 or 3 per group), colored pencils (4 or 5 different colors per group)\",\n  { matchCase: true, matchWholeWord: false }\n);\nawait context.sync();\nif (resources.items.length === 0) {\n  throw new Error(\"Could not find resources text\");\n}\nresources.items[0].insertText(\n  \"Papier vierge (5 par groupe d'\u00e9l\u00e8ves), stylo (2 ou 3 par groupe), crayons de couleur (4 ou 5 couleurs diff\u00e9rentes par groupe)\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 5) In-person session tip cell\nconst tip = body.search(\n  \"In case the session takes place in person, the facilitator can draw some lines on the floor with tape so that some student volunteer can try and walk along them after the introduction to the problem.\",\n  { matchCase: true, matchWholeWord: false }\n);\nawait context.sync();\nif (tip.items.length === 0) {\n  throw new Error(\"Could not find in-person tip text\");\n}\ntip.items[0].insertText(\n  \"Si la session a lieu en personne, l'animateur peut tracer des lignes sur le sol avec du ruban adh\u00e9sif afin qu'un \u00e9l\u00e8ve volontaire puisse essayer de marcher le long de ces lignes apr\u00e8s l'introduction du probl\u00e8me.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Translate the \"Path walking Snail\" facilitator guide table into French.\n# Each edit below grabs a *fresh* Range over the whole document content,\n# runs Find to locate a unique piece of text, and then assigns Range.Text\n# directly (instead of Find.Execute(..., Replace:=...)) so Word's\n# smart-quote AutoCorrect never mangles the straight apostrophes used in\n# the French replacement text. Re-deriving the range from $d.Content each\n# time (rather than reusing/caching a Range or Table.Cell handle across\n# edits) keeps each Find anchored to the live, current document text.\n\n$d = $word.ActiveDocument\n\n# 1) Title row, first run: \"The \" -> \"L'escargot marchant \"\n#    (the bold \"p\" run that follows stays untouched, see edit #2)\n$r1 = $d.Content\n$find1 = $r1.Find\n$find1.Text = \"The \"\n$found1 = $find1.Execute()\nif (-not $found1) {\n  throw \"Could not find title run 'The '\"\n}\n$r1.Text = \"L'escargot marchant \"\n\n# 2) Title row, third run: \"ath walking Snail\" -> \"ath\"\n$r2 = $d.Content\n$find2 = $r2.Find\n$find2.Text = \"ath walking Snail\"\n$found2 = $find2.Execute()\nif (-not $found2) {\n  throw \"Could not find title run 'ath walking Snail'\"\n}\n$r2.Text = \"ath\"\n\n# 3) Objective(s) cell\n$r3 = $d.Content\n$find3 = $r3.Find\n$find3.Text = \"Understand how to model a difficult mathematical problem. Learn how to solve a puzzle by means of a coloring of the plane.\"\n$found3 = $find3.Execute()\nif (-not $found3) {\n  throw \"Could not find objective text\"\n}\n$r3.Text = \" Apprenez \u00e0 r\u00e9soudre une \u00e9nigme \u00e0 l'aide d'un coloriage du plan.\"\n\n# 4) Required resources cell\n$r4 = $d.Content\n$find4 = $r4.Find\n$find4.Text = \"Blank paper (5 per group of students), pen (2 or 3 per group), colored pencils (4 or 5 different colors per group)\"\n$found4 = $find4.Execute()\nif (-not $found4) {\n  throw \"Could not find resources text\"\n}\n$r4.Text = \"Papier vierge (5 par groupe d'\u00e9l\u00e8ves), stylo (2 ou 3 par groupe), crayons de couleur (4 ou 5 couleurs diff\u00e9rentes par groupe)\"\n\n# 5) In-person session tip cell\n$r5 = $d.Content\n$find5 = $r5.Find\n$find5.Text = \"In case the session takes place in person, the facilitator can draw some lines on the floor with tape so that some student volunteer can try and walk along them after the introduction to the problem.\"\n$found5 = $find5.Execute()\nif (-not $found5) {\n  throw \"Could not find in-person tip text\"\n}\n$r5.Text = \"Si la session a lieu en personne, l'animateur peut tracer des lignes sur le sol avec du ruban adh\u00e9sif afin qu'un \u00e9l\u00e8ve volontaire puisse essayer de marcher le long de ces lignes apr\u00e8s l'introduction du probl\u00e8me.\"\n"}
